# Update reports - 2026-01-28 04:11
# A new scan (2026-01-28 04:11:36) found three new publication IDs
# (51276, 51283, 51285) and two of them turned into new report rows
# on the "Reports" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "Reports" sheet/table: append the two newly-published reports
# ---------------------------------------------------------------
$reportsSheet = $wb.Worksheets.Item("Reports")
$reportsTable = $reportsSheet.ListObjects.Item("Reports")

$row26 = $reportsTable.ListRows.Add()
$row27 = $reportsTable.ListRows.Add()

# Row 26 - Public Accounts Committee / 63rd Report
$reportsSheet.Cells.Item(26, 1).Value  = "'51276"
$reportsSheet.Cells.Item(26, 2).Value  = "HC 1239"
$reportsSheet.Cells.Item(26, 3).Value  = "2024-26"
$reportsSheet.Cells.Item(26, 4).Value  = "Public Accounts Committee"
$reportsSheet.Cells.Item(26, 5).Value  = "Commons"
$reportsSheet.Cells.Item(26, 6).Value  = "Increasing police productivity"
$reportsSheet.Cells.Item(26, 7).Value  = "63rd Report"
$reportsSheet.Cells.Item(26, 8).Value  = "'2026-01-28"
$reportsSheet.Cells.Item(26, 9).Value  = "00:01:00"
$reportsSheet.Cells.Item(26, 10).Value = "0:20:31"
$reportsSheet.Cells.Item(26, 11).Value = "4:10:36"

# Row 27 - Transport Committee / 5th Report
$reportsSheet.Cells.Item(27, 1).Value  = "'51283"
$reportsSheet.Cells.Item(27, 2).Value  = "HC 1234"
$reportsSheet.Cells.Item(27, 3).Value  = "2024-26"
$reportsSheet.Cells.Item(27, 4).Value  = "Transport Committee"
$reportsSheet.Cells.Item(27, 5).Value  = "Commons"
$reportsSheet.Cells.Item(27, 6).Value  = "Engine for growth: securing skills for transport manufacturing"
$reportsSheet.Cells.Item(27, 7).Value  = "5th Report"
$reportsSheet.Cells.Item(27, 8).Value  = "'2026-01-28"
$reportsSheet.Cells.Item(27, 9).Value  = "00:01:00"
$reportsSheet.Cells.Item(27, 10).Value = "0:20:31"
$reportsSheet.Cells.Item(27, 11).Value = "4:10:36"

# ---------------------------------------------------------------
# 2. "Scans" sheet/table: log the new scan that found the three IDs
# ---------------------------------------------------------------
$scansSheet = $wb.Worksheets.Item("Scans")
$scansTable = $scansSheet.ListObjects.Item("Scans")

$row8 = $scansTable.ListRows.Add()

$scansSheet.Cells.Item(8, 1).Value = "'2026-01-28"
$scansSheet.Cells.Item(8, 2).Value = "04:11:36"
$scansSheet.Cells.Item(8, 3).Value = "51276, 51283, 51285"
